# Add a new vulnerability definition row ("PTV-NET-IDENT-ACTIVE-MLDNVERDEV2")
# into the defined-vulnerability list, right after the existing
# "PTV-NET-IDENT-ACTIVE-MLDNVERDEV" row (row 9), pushing every following
# row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (shifts rows 10.. down to 11..)
$ws.Rows("10:10").Insert()

# Populate the new row with the new vulnerability definition
$ws.Range("A10").Value = "Node"
$ws.Range("B10").Value = "a,a+"
$ws.Range("C10").Value = "PTV-NET-IDENT-ACTIVE-MLDNVERDEV2"
$ws.Range("D10").Value = "Device only responds to illegitimate MLDv1 queries even though MLDv2 queries are sent, possibly downgraded"

# Update the saved selection/viewport to reflect where the author left off
$ws.Range("H18").Select()
